$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 6 - "Ver dispositivos del cliente" test case: fill in the previously
# empty columns C..L with the actual test-case content.
# ---------------------------------------------------------------------------
$ws.Range("C6").Value = "Positivo"
$ws.Range("D6").Value = "eCenter"
$ws.Range("E6").Value = "Cliente en estado ACTIVO seleccionado"
$ws.Range("F6").Value = "1. Clic en Opciones.`n2. Seleccionar Ver dispositivos."
$ws.Range("G6").Value = "N/A"
$ws.Range("H6").Value = "El sistema abre un modal que lista los dispositivos asociados al cliente."
$ws.Range("I6").Value = "el modal se visualiza correctamente"
$ws.Range("J6").Value = "OK"
$ws.Range("K6").Value = "SI"
$ws.Range("L6").Value = "N/A"

# ---------------------------------------------------------------------------
# Row 7 - brand-new test case CP_GESCLSERDOM_006: view & send documents
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "CP_GESCLSERDOM_006"
$ws.Range("B7").Value = "Ver y enviar documentos (Acta de instalación y Contrato)"
$ws.Range("C7").Value = "Positivo"
$ws.Range("D7").Value = "eCenter"
$ws.Range("E7").Value = "Cliente con `ndocumentos asociados"
$ws.Range("F7").Value = "1. Clic en Opciones.`n2. Seleccionar Ver documentos.`n3. Para cada documento:`n a. Clic Ver documento.`n b. Clic Enviar al correo.`n c. Clic Descargar."
$ws.Range("G7").Value = "N/A"
$ws.Range("H7").Value = "El sistema abre un modal para visualizar los documentos del cliente(Acta de instalación y Contrato)"
$ws.Range("I7").Value = "el modal se visualiza correctamente"
$ws.Range("J7").Value = "OK"
$ws.Range("K7").Value = "SI"
$ws.Range("L7").Value = "N/A"

# ---------------------------------------------------------------------------
# Formatting touch-ups on the two rows so the new cells read like the rest
# of the table (vertical-centered, wrapped text, matching borders/fonts).
# ---------------------------------------------------------------------------
$ws.Range("C6:I6").VerticalAlignment = -4108
$ws.Range("C6:I6").WrapText = $true
$ws.Range("J5").Copy()
$ws.Range("J6").PasteSpecial(-4122)
$ws.Range("K5").Copy()
$ws.Range("K6").PasteSpecial(-4122)
$ws.Range("L6").VerticalAlignment = -4108
$ws.Range("L6").WrapText = $true

$ws.Range("A7:B7").VerticalAlignment = -4108
$ws.Range("A7:B7").WrapText = $true
$ws.Range("C7").VerticalAlignment = -4108
$ws.Range("C7").WrapText = $false
$ws.Range("D7").VerticalAlignment = -4108
$ws.Range("D7").WrapText = $false
$ws.Range("E7").VerticalAlignment = -4108
$ws.Range("E7").WrapText = $true
$ws.Range("F7").VerticalAlignment = -4108
$ws.Range("F7").WrapText = $true
$ws.Range("G7").VerticalAlignment = -4108
$ws.Range("G7").WrapText = $true
$ws.Range("H7").VerticalAlignment = -4108
$ws.Range("H7").WrapText = $true
$ws.Range("I2").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("I7").Value = "el modal se visualiza correctamente"
$ws.Range("J5").Copy()
$ws.Range("J7").PasteSpecial(-4122)
$ws.Range("K5").Copy()
$ws.Range("K7").PasteSpecial(-4122)
$ws.Range("L7").VerticalAlignment = -4108
$ws.Range("L7").WrapText = $true

# ---------------------------------------------------------------------------
# Selection matches where the author was working when the test cases were
# added (scrolled toward the new rows, cursor left on the last cell typed).
# ---------------------------------------------------------------------------
$ws.Range("L6").Select()
